$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Roll the yearly headers forward by one year (1396..1400 -> 1397..1401) ---
$ws.Range("E8").Value  = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F8").Value  = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G8").Value  = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H8").Value  = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I8").Value  = "دوازده ماهه منتهی به 1401/12"

$ws.Range("E24").Value = "دوازده ماهه منتهی به 1397/12"
$ws.Range("F24").Value = "دوازده ماهه منتهی به 1398/12"
$ws.Range("G24").Value = "دوازده ماهه منتهی به 1399/12"
$ws.Range("H24").Value = "دوازده ماهه منتهی به 1400/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# --- Update "هزینه های عمومی و اداری" detail table (rows 10-20) ---
# هزینه حمل و نقل و انتقال
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 20
$ws.Range("G10").Value = 41
$ws.Range("H10").Value = 23
$ws.Range("I10").Value = 252

# هزینه خدمات پس از فروش (row 11) - unchanged (all zero)

# حق العمل و کمیسیون فروش (row 12) - unchanged (all zero)

# هزینه تبلیغات
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 666
$ws.Range("I13").Value = 920

# هزینه مواد مصرفی
$ws.Range("E14").Value = 49
$ws.Range("F14").Value = 96
$ws.Range("G14").Value = 22
$ws.Range("H14").Value = 189
$ws.Range("I14").Value = 1

# هزینه انرژی (آب، برق، گاز و سوخت)
$ws.Range("E15").Value = 120
$ws.Range("F15").Value = 103
$ws.Range("G15").Value = 136
$ws.Range("H15").Value = 563
$ws.Range("I15").Value = 1898

# هزینه استهلاک
$ws.Range("E16").Value = 835
$ws.Range("F16").Value = 1020
$ws.Range("G16").Value = 1471
$ws.Range("H16").Value = 3177
$ws.Range("I16").Value = 2994

# هزینه حقوق و دستمزد
$ws.Range("E17").Value = 25021
$ws.Range("F17").Value = 36849
$ws.Range("G17").Value = 51439
$ws.Range("H17").Value = 76608
$ws.Range("I17").Value = 113371

# هزینه مطالبات مشکوک الوصول (row 18) - unchanged (all zero)

# سایر هزینه ها
$ws.Range("E19").Value = 25597
$ws.Range("F19").Value = 33068
$ws.Range("G19").Value = 63527
$ws.Range("H19").Value = 96681
$ws.Range("I19").Value = 94433

# جمع (total)
$ws.Range("E20").Value = 51624
$ws.Range("F20").Value = 71156
$ws.Range("G20").Value = 116636
$ws.Range("H20").Value = 177907
$ws.Range("I20").Value = 213869

# --- Update "تعداد پرسنل" (personnel count) table (rows 26-27) ---
# تعداد پرسنل غیر تولیدی شرکت
$ws.Range("E26").Value = 145
$ws.Range("F26").Value = 159
$ws.Range("G26").Value = 160
$ws.Range("H26").Value = 151
$ws.Range("I26").Value = 148

# تعداد پرسنل تولیدی شرکت
$ws.Range("E27").Value = 85
$ws.Range("F27").Value = 90
$ws.Range("G27").Value = 83
$ws.Range("H27").Value = 102
$ws.Range("I27").Value = 102
